# WAT new test cases
# Adds 5 new worksheet rows (63-67) describing new WAT-928..WAT-932 publication
# recommendation test cases, matching the shape/values/styling used by the
# existing rows in the WoS_AuthorTransformation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WoS_AuthorTransformation")

# ---------------------------------------------------------------------------
# 1. New row data (values only - columns A..L, same layout as existing rows)
# ---------------------------------------------------------------------------

# Row 63
$ws.Range("A63").Value = "WAT-928"
$ws.Range("B63").Value = "Verify that publication recommendations display on author record using authorId"
$ws.Range("C63").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D63").Value = "/author/recommend/publications"
$ws.Range("E63").Value = "GET"
$ws.Range("G63").Value = "?authorId=7781649"
$ws.Range("J63").Value = "status=200"
$ws.Range("K63").Value = "hits[0].authorId||hits[0].ut||hits[0].title||hits[0].authors||hits[0].journal||hits[0].volume||hits[0].issue||hits[0].published||hits[0].page||hits[0].timesCited||hits[0].score"

# Row 64
$ws.Range("A64").Value = "WAT-929"
$ws.Range("B64").Value = "Verify that max 3 publication recommendations display on author record"
$ws.Range("C64").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D64").Value = "/author/recommend/publications"
$ws.Range("E64").Value = "GET"
$ws.Range("G64").Value = "?authorId=7781649"
$ws.Range("J64").Value = "status=200"
$ws.Range("K64").Value = "hits[0].authorId||hits[1].authorId||hits[2].authorId"

# Row 65
$ws.Range("A65").Value = "WAT-930"
$ws.Range("B65").Value = "Verify that publication recommendations display on author record using authorId and name"
$ws.Range("C65").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D65").Value = "/author/recommend/publications"
$ws.Range("E65").Value = "GET"
$ws.Range("G65").Value = "?authorId=7781649&name=FABBRI, F."
$ws.Range("J65").Value = "status=200"
$ws.Range("K65").Value = "hits[0].authorId||hits[0].ut||hits[0].title||hits[0].authors||hits[0].journal||hits[0].volume||hits[0].issue||hits[0].published||hits[0].page||hits[0].timesCited||hits[0].score"

# Row 66
$ws.Range("A66").Value = "WAT-931"
$ws.Range("B66").Value = "Verify that publication recommendations should not display for missing authorId"
$ws.Range("C66").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D66").Value = "/author/recommend/publications"
$ws.Range("E66").Value = "GET"
$ws.Range("G66").Value = "?name=FABBRI, F."
$ws.Range("J66").Value = "status=400"

# Row 67
$ws.Range("A67").Value = "WAT-932"
$ws.Range("B67").Value = "Verify that each publication recommendations should have score"
$ws.Range("C67").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D67").Value = "/author/recommend/publications"
$ws.Range("E67").Value = "GET"
$ws.Range("G67").Value = "?authorId=7781649&name=FABBRI, F."
$ws.Range("J67").Value = "status=200"
$ws.Range("K67").Value = "hits[0].score||hits[1].score||hits[2].score"

Write-Host "values written"
